$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update status ("Status" column, I) for several existing stories to "Done"
# ---------------------------------------------------------------------------
$doneRows = 2, 8, 9, 10, 11, 12, 15, 16
foreach ($r in $doneRows) {
    $ws.Cells.Item($r, 9).Value = "Done"
}

# ---------------------------------------------------------------------------
# 2) Widen column D (the "I want to" column)
# ---------------------------------------------------------------------------
$ws.Columns.Item(4).ColumnWidth = 66 - (5 / 6)

# ---------------------------------------------------------------------------
# 3) Append three new user stories as rows 21-23
# ---------------------------------------------------------------------------
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "Functional"
$ws.Cells.Item(21, 3).Value = "Client"
$ws.Cells.Item(21, 4).Value = "see markers with different colours based on price"
$ws.Cells.Item(21, 5).Value = "I can visually understand the price of different locations"
$ws.Cells.Item(21, 7).Value = 2
$ws.Cells.Item(21, 8).Value = 2
$ws.Cells.Item(21, 9).Value = "In Production"

$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "Functional"
$ws.Cells.Item(22, 3).Value = "Client"
$ws.Cells.Item(22, 4).Value = "see a visual highlight of the corresponding facility when I click on a marker "
$ws.Cells.Item(22, 5).Value = "I can read what cell in the table corresponds to the marker"
$ws.Cells.Item(22, 7).Value = 3
$ws.Cells.Item(22, 8).Value = 3
$ws.Cells.Item(22, 9).Value = "In Production"

$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "Functional"
$ws.Cells.Item(23, 3).Value = "Client"
$ws.Cells.Item(23, 4).Value = "Click on a facility name and have it behave like clicking a marker"
$ws.Cells.Item(23, 5).Value = "I can navigate the website easier"
$ws.Cells.Item(23, 7).Value = 3
$ws.Cells.Item(23, 8).Value = 3
$ws.Cells.Item(23, 9).Value = "In Production"

# ---------------------------------------------------------------------------
# 4) Grow Table1 so the new rows belong to it, and keep the autofilter in sync
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:I24"))

# ---------------------------------------------------------------------------
# 5) Match the saved selection / active cell
# ---------------------------------------------------------------------------
$ws.Range("I21").Select()

# ---------------------------------------------------------------------------
# 6) Reflect the window size recorded in the workbook view
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.Width = 28800
$win.Height = 12300
